# "A mistake in the numbers alerted by Tim."
# Two data-entry typos on the "Full Titan" sheet are corrected:
#   B5: 265  -> 256   (message size, should be a power of two)
#   B9: 4069 -> 4096  (message size, should be a power of two)
# All the dependent formula cells (I/J/K/L/M columns) recompute
# automatically, and the scatter charts that plot column B against
# those derived columns pick up the corrected values on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full Titan")

$ws.Range("B5").Value = 256
$ws.Range("B9").Value = 4096
